$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Revert the "datetimeFigureOut" date placeholder text on the slide
#    master and every slide layout: 13.01.2025 -> 10.01.2025
# ---------------------------------------------------------------------------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.TextRange.Text -eq "13.01.2025") {
            $shape.TextFrame.TextRange.Text = "10.01.2025"
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shape = $layout.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq "13.01.2025") {
                $shape.TextFrame.TextRange.Text = "10.01.2025"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Revert the wording on the Slide 3 title placeholder:
#      "Static Class "            -> "Static Class static constructor "
#      "Registered"               -> "injected"
#      (and drop the soft line break that used to separate
#       "...Services " from "@inject ...")
#
#    Edits are applied as targeted Characters(start, length) replacements
#    (right-to-left, so earlier offsets stay valid) instead of rewriting
#    the whole TextRange, which keeps the untouched line breaks / runs of
#    the paragraph intact.
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$title = $slide3.Shapes.Item(1)
$tr = $title.TextFrame.TextRange

# "Services \vinject" -> "Services @inject" (removes the break before "@inject")
$tr.Characters(77, 11).Text = "Services @"

# "Registered" -> "injected"
$tr.Characters(66, 10).Text = "injected"

# "Static Class " -> "Static Class static constructor "
$tr.Characters(29, 13).Text = "Static Class static constructor "
